$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the pin-header row (WJ15EDGVC-3.5-3P) from the button panel BOM.
# The row itself stays in place (dimension is unchanged); only its contents
# are cleared, same as selecting A5:F5 and pressing Delete in Excel.
$ws.Range("A5:F5").ClearContents()
$ws.Range("A5:F5").Select()
